$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.572.87"
$ws.Range("E2").Value = "  +5.09%  "
$ws.Range("D3").Value = "3.106.61"
$ws.Range("E3").Value = "  +3.86%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.34"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.07"
$ws.Range("E6").Value = "  +9.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "3.105.19"
$ws.Range("E8").Value = "  +4.03%  "
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.08"
$ws.Range("E10").Value = "  +17.94%  "
$ws.Range("E11").Value = "  +5.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("E12").Value = "  +4.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("E13").Value = "  +4.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.39"
$ws.Range("E14").Value = "  +4.15%  "
$ws.Range("D15").Value = "3.603.55"
$ws.Range("E15").Value = "  +3.87%  "
$ws.Range("D16").Value = "64.619.86"
$ws.Range("E16").Value = "  +4.98%  "
$ws.Range("D17").Value = "3.109.31"
$ws.Range("E17").Value = "  +4.05%  "
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.78"
$ws.Range("E19").Value = "  +2.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "483.55"
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.81"
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.63"
$ws.Range("E22").Value = "  +9.40%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.675"
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.28"
$ws.Range("E24").Value = "  +10.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.02"
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  +4.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.08"
$ws.Range("E28").Value = "  +5.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.08"
$ws.Range("E29").Value = "  +8.77%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.16"
$ws.Range("E31").Value = "  +2.71%  "
$ws.Range("E32").Value = "  +3.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.45"
$ws.Range("E33").Value = "  +6.20%  "
$ws.Range("E34").Value = "  +5.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.28"
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.14"
$ws.Range("E36").Value = "  +4.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "465.81"
$ws.Range("E37").Value = "  +4.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0409"
$ws.Range("E38").Value = "  +7.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0828"
$ws.Range("E39").Value = "  +4.68%  "
$ws.Range("D40").Value = "3.033.66"
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.117"
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.30"
$ws.Range("E42").Value = "  +2.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.70"
$ws.Range("E43").Value = "  +14.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "28.25"
$ws.Range("E44").Value = "  +10.67%  "
$ws.Range("E45").Value = "  +8.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.11"
$ws.Range("E47").Value = "  +8.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.113"
$ws.Range("E48").Value = "  +4.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "119.13"
$ws.Range("E49").Value = "  +4.21%  "
$ws.Range("D50").Value = "0.0₃0517"
$ws.Range("E50").Value = "  +6.94%  "
$ws.Range("E51").Value = "  +3.38%  "
